$wb = $excel.ActiveWorkbook

# --- tc002: selection changes (used to be the active/selected tab before the run) ---
$tc002 = $wb.Worksheets.Item("tc002")
$tc002.Activate()
$tc002.Range("A1:A2").Select()

# --- tc004: selection C8 -> C14 ---
$tc004 = $wb.Worksheets.Item("tc004")
$tc004.Activate()
$tc004.Range("C14").Select()

# --- tc005: selection B1:B2 -> D18 ---
$tc005 = $wb.Worksheets.Item("tc005")
$tc005.Activate()
$tc005.Range("D18").Select()

# --- tc006: selection A2 -> A11 ---
$tc006 = $wb.Worksheets.Item("tc006")
$tc006.Activate()
$tc006.Range("A11").Select()

# --- Add the new tc003 sheet after the last sheet (tc019) - ran the suite and report generated ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tc003 = $wb.Worksheets.Add($null, $lastSheet)
$tc003.Name = "tc003"

$tc003.Range("A1").Value = "projectName"
$tc003.Range("A2").Value = "STG- PulseCodeOnAzureCloud"

# Match formatting used for the same value elsewhere in the workbook (explicit-black font style)
$tc002.Range("A2").Copy()
$tc003.Range("A2").PasteSpecial(-4122)

# Make tc003 the active sheet/tab with its own selection, as the newest report
$tc003.Activate()
$tc003.Range("D11").Select()
